# Add a new 'removeJob' automation-script entry to the testData sheet,
# mirroring the existing "removeCompany" block (rows 20-21) as rows 32-33:
#   row 32 (header/label row) -> TrainScheduling_ltrailways_removeJob / CompanyManagement.removeJob
#   row 33 (parameter row)    -> rowIndex / isRemove values for that call

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testData")

# --- Copy formatting only (fills/borders/number-format) from the
# "removeCompany" template block into the new block, restricted to the
# exact columns that carry content so we don't materialize extra cells.
$ws.Range("A20:E20").Copy()
$ws.Range("A32:E32").PasteSpecial(-4122) | Out-Null

$ws.Range("A21:E21").Copy()
$ws.Range("A33:E33").PasteSpecial(-4122) | Out-Null

$ws.Range("G21").Copy()
$ws.Range("G33").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Match the template's row height for the two new rows.
$ws.Rows("32").RowHeight = 15.95
$ws.Rows("33").RowHeight = 15.95

# Helper: force a literal text value (avoids "true"/"false"/"1" style
# strings being auto-coerced to Boolean/Number types by plain .Value
# assignment) by round-tripping through a formula + paste-values.
function Set-TextValue($range, $text) {
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = "=""" + $text + """"
    $scratch.Copy()
    $range.PasteSpecial(-4163) | Out-Null
    $scratch.Clear() | Out-Null
    $excel.CutCopyMode = $false
}

# --- Row 32: header/label row for the removeJob call.
$ws.Range("A32").Value = "TrainScheduling_ltrailways_removeJob"
$ws.Range("B32").Value = "1"
$ws.Range("C32").Value = "CompanyManagement.removeJob"
$ws.Range("D32").Value = "rowIndex"
$ws.Range("E32").Value = "isRemove"

# --- Row 33: parameter values row for the removeJob call.
$ws.Range("A33").Value = "TrainScheduling_ltrailways_removeJob"
$ws.Range("B33").Value = "1"
$ws.Range("C33").Value = "CompanyManagement.removeJob"
$ws.Range("D33").Value = "1"
Set-TextValue $ws.Range("E33") "true"

$ws.Range("E33").Select()
